$wb = $excel.ActiveWorkbook

# --- Update summary text on "Hoja1"!A1 with the new conversion rates ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.72 = 6389.72 pesos`n✅ 6389.72 pesos = 1.71 = 907.7 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Update the rate figures on "tasas" sheet ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 580
$wsTasas.Range("O10").Value = 3706.04

$wsTasas.Range("N12").Value = 3730.93
$wsTasas.Range("O12").Value = 530
